$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextValue 'D2' '28.997.37'
Set-TextValue 'E2' '  +0.44%  '
Set-TextValue 'D3' '1.919.18'
Set-TextValue 'E3' '  +1.63%  '
Set-TextValue 'D4' '1.003'
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '324.96'
Set-TextValue 'E5' '  +0.48%  '
Set-TextValue 'D6' '1.002'
Set-TextValue 'E6' '  +0.03%  '
Set-TextValue 'D7' '0.4588'
Set-TextValue 'E7' '  +0.09%  '
Set-TextValue 'D8' '0.3811'
Set-TextValue 'E8' '  +0.07%  '
Set-TextValue 'D9' '0.07737'
Set-TextValue 'E9' '  +0.29%  '
Set-TextValue 'D10' '0.9767'
Set-TextValue 'E10' '  +1.18%  '
Set-TextValue 'D11' '22.72'
Set-TextValue 'E11' '  +3.19%  '
Set-TextValue 'D12' '1.950.27'
Set-TextValue 'E12' '  +3.14%  '
Set-TextValue 'D13' '5.688'
Set-TextValue 'E13' '  +0.46%  '
Set-TextValue 'E14' '  +0.28%  '
Set-TextValue 'D15' '0.07012'
Set-TextValue 'E15' '  -0.31%  '
Set-TextValue 'D16' '1.005'
Set-TextValue 'E16' '  +0.11%  '
Set-TextValue 'D17' '84.16'
Set-TextValue 'E17' '  +1.11%  '
Set-TextValue 'D18' '0.000009491'
Set-TextValue 'E18' '  -0.07%  '
Set-TextValue 'D19' '16.68'
Set-TextValue 'D20' '1.002'
Set-TextValue 'E20' '  +0.06%  '
Set-TextValue 'D21' '29.005.51'
Set-TextValue 'E21' '  +0.59%  '
Set-TextValue 'D22' '5.341'
Set-TextValue 'E22' '  +0.98%  '
Set-TextValue 'E23' '  +1.04%  '
Set-TextValue 'B24' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C24' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D24' '2.125.95'
Set-TextValue 'E24' '  +0.12%  '
Set-TextValue 'B25' 'Toncoin'
Set-TextValue 'C25' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D25' '2.068'
Set-TextValue 'E25' '  -0.48%  '
Set-TextValue 'B26' 'Monero'
Set-TextValue 'C26' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D26' '157.22'
Set-TextValue 'E26' '  +0.56%  '
Set-TextValue 'B27' 'EthereumClassic'
Set-TextValue 'C27' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D27' '19.04'
Set-TextValue 'E27' '  +0.36%  '
Set-TextValue 'B28' 'InternetComputer(DFINITY)'
Set-TextValue 'C28' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D28' '5.643'
Set-TextValue 'E28' '  +1.18%  '
Set-TextValue 'B29' 'BitcoinCash'
Set-TextValue 'C29' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D29' '117.64'
Set-TextValue 'E29' '  +0.52%  '
Set-TextValue 'B30' 'LidoDAOToken'
Set-TextValue 'C30' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D30' '1.830'
Set-TextValue 'E30' '  +1.16%  '
Set-TextValue 'B31' 'Stellar'
Set-TextValue 'C31' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D31' '0.09321'
Set-TextValue 'E31' '  +0.76%  '
Set-TextValue 'B32' 'ImmutableX'
Set-TextValue 'C32' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D32' '0.8592'
Set-TextValue 'E32' '  +1.32%  '
Set-TextValue 'B33' 'Filecoin'
Set-TextValue 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '5.087'
Set-TextValue 'E33' '  +0.38%  '
Set-TextValue 'B34' 'ARBITRUM'
Set-TextValue 'C34' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D34' '1.239'
Set-TextValue 'E34' '  +1.18%  '
Set-TextValue 'B35' 'HuobiToken'
Set-TextValue 'C35' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D35' '3.011'
Set-TextValue 'E35' '  -0.09%  '
Set-TextValue 'B36' 'TrustWalletToken'
Set-TextValue 'C36' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D36' '1.158'
Set-TextValue 'E36' '  +1.27%  '
Set-TextValue 'B37' 'Hedera'
Set-TextValue 'C37' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D37' '0.05671'
Set-TextValue 'E37' '  -0.01%  '
Set-TextValue 'B38' 'Frax'
Set-TextValue 'C38' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D38' '1.001'
Set-TextValue 'E38' '  -0.06%  '
Set-TextValue 'B39' 'VeChain'
Set-TextValue 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D39' '0.02044'
Set-TextValue 'E39' '  +0.56%  '
Set-TextValue 'B40' 'MXToken'
Set-TextValue 'C40' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D40' '3.096'
Set-TextValue 'E40' '  +14.84%  '
Set-TextValue 'B41' 'FraxShare'
Set-TextValue 'C41' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D41' '7.422'
Set-TextValue 'E41' '  +0.37%  '
Set-TextValue 'B42' 'TheSandbox'
Set-TextValue 'C42' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D42' '0.5493'
Set-TextValue 'E42' '  +0.24%  '
Set-TextValue 'B43' 'Algorand'
Set-TextValue 'C43' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D43' '0.1751'
Set-TextValue 'E43' '  +0.28%  '
Set-TextValue 'B44' 'Aptos'
Set-TextValue 'C44' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D44' '9.354'
Set-TextValue 'E44' '  +2.16%  '
Set-TextValue 'D45' '0.000002843'
Set-TextValue 'E45' '  -1.59%  '
Set-TextValue 'B46' 'RenderToken'
Set-TextValue 'C46' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D46' '2.186'
Set-TextValue 'E46' '  +5.47%  '
Set-TextValue 'B47' 'Decentraland'
Set-TextValue 'C47' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D47' '0.5172'
Set-TextValue 'E47' '  +0.22%  '
Set-TextValue 'B48' 'EnergySwap'
Set-TextValue 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '11.26'
Set-TextValue 'E48' '  +0.21%  '
Set-TextValue 'B49' 'Cronos'
Set-TextValue 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D49' '0.06899'
Set-TextValue 'E49' '  +1.61%  '
Set-TextValue 'B50' 'Quant'
Set-TextValue 'C50' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D50' '110.18'
Set-TextValue 'E50' '  -0.96%  '
Set-TextValue 'B51' 'NEARProtocol'
Set-TextValue 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D51' '1.761'
Set-TextValue 'E51' '  -0.89%  '
